$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (H1) onto the new
# header cells I1 and J1, so they match the bold/bordered/centered style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set header labels for the new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for column I (I0) and J (IF), rows 2-7
$values = @(
    @(7, 8),
    @(7, 7),
    @(9, 9),
    @(6, 7),
    @(3, 3),
    @(3, 3)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
